$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Prep: copy the formatting of the last existing data row (row 201) down
# onto the new rows 202-207 so the cell styles match (date / hyperlink /
# plain-text / remark styles). Only touch the columns that actually hold a
# cell in each new row (row 204 has no E value, rows 205-207 have no N
# value at all).
# ---------------------------------------------------------------------------
$srcRow = 201
$newRows = 202,203,204,205,206,207
$rowCols = @{
    "202" = @("A","B","C","D","E","N")
    "203" = @("A","B","C","D","E","N")
    "204" = @("A","B","C","D","N")
    "205" = @("A","B","C","D","E")
    "206" = @("A","B","C","D","E")
    "207" = @("A","B","C","D","E")
}
foreach ($r in $newRows) {
    foreach ($col in $rowCols["$r"]) {
        $ws.Range("$col$srcRow").Copy() | Out-Null
        $ws.Range("$col$r").PasteSpecial(-4122) | Out-Null
    }
}

# ---------------------------------------------------------------------------
# Fill in the non-text (numeric) values first - these never touch the
# shared-string table so their order doesn't matter.
# ---------------------------------------------------------------------------
$ws.Range("A202").Value = 43851
$ws.Range("B202").Value = 1489

$ws.Range("A203").Value = 44218
$ws.Range("B203").Value = 989

$ws.Range("A204").Value = 44222
$ws.Range("B204").Value = 1128

$ws.Range("A205").Value = 44224
$ws.Range("B205").Value = 724

$ws.Range("A206").Value = 44249
$ws.Range("B206").Value = 766

$ws.Range("A207").Value = 44250
$ws.Range("B207").Value = 1052

# ---------------------------------------------------------------------------
# Now add the brand-new text values, in the exact order the author typed
# them, so new shared-string entries come out in the same sequence as the
# target workbook (排名：40931, 找到最小生成树..., 图，最小生成树，并查集,
# 寻找数组的中心索引, 托普利茨矩阵, 数组形式的整数加法, 爱生气的书店老板,
# 等价多米诺骨牌的数量).
# ---------------------------------------------------------------------------
$ws.Range("N202").Value = "排名：40931"
$ws.Range("C202").Value = "找到最小生成树里的关键边和伪关键边"
$ws.Range("E200").Value = "图，最小生成树，并查集"
$ws.Range("C205").Value = "寻找数组的中心索引"
$ws.Range("C206").Value = "托普利茨矩阵"
$ws.Range("C203").Value = "数组形式的整数加法"
$ws.Range("C207").Value = "爱生气的书店老板"
$ws.Range("C204").Value = "等价多米诺骨牌的数量"

# ---------------------------------------------------------------------------
# Remaining D/E values - these all reuse pre-existing shared strings
# (简单/中等/困难/数组/贪心) so order is irrelevant here.
# ---------------------------------------------------------------------------
$ws.Range("D202").Value = "困难"
$ws.Range("E202").Value = "图，最小生成树，并查集"

$ws.Range("D203").Value = "简单"
$ws.Range("E203").Value = "数组"

$ws.Range("D204").Value = "简单"

$ws.Range("D205").Value = "简单"
$ws.Range("E205").Value = "数组"

$ws.Range("D206").Value = "简单"
$ws.Range("E206").Value = "数组"

$ws.Range("D207").Value = "中等"
$ws.Range("E207").Value = "贪心"

# ---------------------------------------------------------------------------
# Hyperlinks for the problem numbers in column B (added in the order the
# author committed them: 1489 first, then 724, 766, 989, 1052, 1128).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B202"), "https://leetcode-cn.com/problems/find-critical-and-pseudo-critical-edges-in-minimum-spanning-tree/", "", "", "1489") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B205"), "https://leetcode-cn.com/problems/find-pivot-index/", "", "", "724") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B206"), "https://leetcode-cn.com/problems/toeplitz-matrix/", "", "", "766") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B203"), "https://leetcode-cn.com/problems/add-to-array-form-of-integer/", "", "", "989") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B207"), "https://leetcode-cn.com/problems/grumpy-bookstore-owner/", "", "", "1052") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B204"), "https://leetcode-cn.com/problems/number-of-equivalent-domino-pairs/", "", "", "1128") | Out-Null

# re-stamp the hyperlink cell style (Hyperlinks.Add recolors the cell) so it
# matches the rest of column B.
foreach ($r in $newRows) {
    $ws.Range("B$srcRow").Copy() | Out-Null
    $ws.Range("B$r").PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# Column E grew wider to fit the longer new tag text.
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 24.5

# ---------------------------------------------------------------------------
# Leave the selection on C206, matching where the author ended up.
# ---------------------------------------------------------------------------
$ws.Range("C206").Select() | Out-Null
